# Review of insurance contract section: restructure Sheet1's header/data
# columns (new fields inserted, some renamed/reordered) and refresh
# shared-string content accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Formatting pass (done first, while the known-good source cells for
#    each style still hold their original formatting) -- copy the
#    "text" (numFmtId 49) style onto every cell that needs it in the
#    final layout, and reset K1 (text style in the original layout)
#    back to the plain/default style it needs going forward.
# ---------------------------------------------------------------------

# K1 currently carries the "text" style but the new layout needs it
# plain -- borrow the plain formatting from A1 (never styled).
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Cells that must end up with the "text" style (numFmtId 49) -- grab it
# from C2, which already carries that style in the original workbook.
$ws.Range("C2").Copy()
$textTargets = "F2","I2","J2","K2","L2","M2","N2","O2","P1","Q1","S2"
foreach ($addr in $textTargets) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2) Content pass -- (re)write every cell's value for the new layout.
# ---------------------------------------------------------------------

$ws.Range("A1").Value = "شناسه قرارداد*"
$ws.Range("B1").Value = "کد الکترونیکی"
$ws.Range("C1").Value = "نام *"
$ws.Range("D1").Value = "نام خانوادگی *"
$ws.Range("E1").Value = "نام پدر"
$ws.Range("F1").Value = "تاریخ تولد"
$ws.Range("G1").Value = "شماره شناسنامه"
$ws.Range("H1").Value = "کد ملی *"
$ws.Range("I1").Value = "تاریخ استخدام"
$ws.Range("J1").Value = "کد جنسیت"
$ws.Range("K1").Value = "وضعیت تاهل"
$ws.Range("L1").Value = "شناسه زیرگروه"
$ws.Range("M1").Value = "شناسه نوع بیمه"
$ws.Range("N1").Value = "شماره دفترچه بیمه پایه"
$ws.Range("O1").Value = "کد بانک"
$ws.Range("P1").Value = "شماره حساب"
$ws.Range("Q1").Value = "تلفن ثابت"
$ws.Range("R1").Value = "شماره همراه *"
$ws.Range("S1").Value = "شماره شبا"
$ws.Range("T1").Value = "وضعیت"

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "تست"
$ws.Range("D2").Value = "تست زاده"
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "1365/01/01"
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = "4342340518"
$ws.Range("I2").Value = "1365/01/01"
$ws.Range("J2").Value = "0 = namoshakhas, 1 = mard, 2 = zan"
$ws.Range("K2").Value = "0 = namoshakhas, 1 = mojarad, 2 = motahel, 3 = motalegheh, 4 = moayal"
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").Value = "True/False"

# ---------------------------------------------------------------------
# 3) Column widths -- best-effort match of the refreshed layout.
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 11.71
$ws.Columns.Item(2).ColumnWidth = 10.29
$ws.Columns.Item(3).ColumnWidth = 10.71
$ws.Columns.Item(4).ColumnWidth = 12.71
$ws.Columns.Item(5).ColumnWidth = 12.71
$ws.Columns.Item(6).ColumnWidth = 10.71
$ws.Columns.Item(7).ColumnWidth = 12.29
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 11.57
$ws.Columns.Item(10).ColumnWidth = 31.86
$ws.Columns.Item(11).ColumnWidth = 65
$ws.Columns.Item(12).ColumnWidth = 22
$ws.Columns.Item(13).ColumnWidth = 22
$ws.Columns.Item(14).ColumnWidth = 22
$ws.Columns.Item(15).ColumnWidth = 22
$ws.Columns.Item(16).ColumnWidth = 16.14
$ws.Columns.Item(17).ColumnWidth = 16.14
$ws.Columns.Item(18).ColumnWidth = 14.86
$ws.Columns.Item(19).ColumnWidth = 12
$ws.Columns.Item(20).ColumnWidth = 10.43

# ---------------------------------------------------------------------
# 4) View state -- selection moves to S2, scrolled so column I is first
#    visible column.
# ---------------------------------------------------------------------

$ws.Range("S2").Select()
